$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.3
$ws.Range("H2").Value = 2.08
$ws.Range("J2").Value = 3.5
$ws.Range("K2").Value = 3.65
$ws.Range("N2").Value = 3.5
$ws.Range("P2").Value = 1.82
$ws.Range("Q2").Value = 2.14
$ws.Range("S2").Value = 3.9
$ws.Range("W2").Value = 1.31
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 8.6
$ws.Range("AA2").Value = 30
$ws.Range("AB2").Value = 14
$ws.Range("AE2").Value = 29
$ws.Range("AF2").Value = 36
$ws.Range("AG2").Value = 17.5
$ws.Range("AH2").Value = 22
$ws.Range("AK2").Value = 70
$ws.Range("AL2").Value = 80
$ws.Range("AN2").Value = 80
$ws.Range("AO2").Value = 19

# Row 3
$ws.Range("F3").Value = 1.46
$ws.Range("G3").Value = 1.57
$ws.Range("W3").Value = 2.74

# Row 4
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2.02
$ws.Range("I4").Value = 4.3
$ws.Range("Q4").Value = 1.96
$ws.Range("R4").Value = 1.38
$ws.Range("S4").Value = 3.45
$ws.Range("V4").Value = 1.31
$ws.Range("W4").Value = 1.99
$ws.Range("Z4").Value = 30
$ws.Range("AD4").Value = 17

# Row 5
$ws.Range("G5").Value = 1.35
$ws.Range("K5").Value = 5.7
$ws.Range("S5").Value = 3.6
$ws.Range("U5").Value = 1.63
$ws.Range("W5").Value = 3.85
$ws.Range("AD5").Value = 50

# Row 6
$ws.Range("F6").Value = 2.7
$ws.Range("H6").Value = 2.52
$ws.Range("I6").Value = 2.8
$ws.Range("K6").Value = 3.85
$ws.Range("N6").Value = 3.8
$ws.Range("S6").Value = 3.1
$ws.Range("V6").Value = 1.55
$ws.Range("AF6").Value = 24

# Row 7
$ws.Range("G7").Value = 1.66
$ws.Range("V7").Value = 1.19
$ws.Range("AB7").Value = 10

# Row 8
$ws.Range("F8").Value = 1.49
$ws.Range("J8").Value = 4.7

# Row 9
$ws.Range("G9").Value = 2.24
$ws.Range("K9").Value = 4.1
$ws.Range("L9").Value = 1.25
$ws.Range("P9").Value = 2.44

# Row 10
$ws.Range("F10").Value = 1.82
$ws.Range("G10").Value = 1.83
$ws.Range("H10").Value = 4.7
$ws.Range("I10").Value = 4.8
$ws.Range("K10").Value = 4.3
$ws.Range("N10").Value = 5.3
$ws.Range("Q10").Value = 1.64
$ws.Range("R10").Value = 1.59
$ws.Range("S10").Value = 2.58
$ws.Range("T10").Value = 1.64
$ws.Range("U10").Value = 2.48
$ws.Range("V10").Value = 1.26
$ws.Range("W10").Value = 2.2
$ws.Range("Y10").Value = 23
$ws.Range("Z10").Value = 38
$ws.Range("AA10").Value = 95
$ws.Range("AK10").Value = 16
$ws.Range("AO10").Value = 38

# Row 11
$ws.Range("F11").Value = 1.9
$ws.Range("G11").Value = 1.92
$ws.Range("H11").Value = 4.2
$ws.Range("N11").Value = 5.3
$ws.Range("O11").Value = 1.2
$ws.Range("R11").Value = 1.6
$ws.Range("S11").Value = 2.56
$ws.Range("W11").Value = 2.08
$ws.Range("AA11").Value = 80
$ws.Range("AI11").Value = 42

# Row 12
$ws.Range("G12").Value = 13
$ws.Range("J12").Value = 7.8
$ws.Range("K12").Value = 8
$ws.Range("O12").Value = 1.08
$ws.Range("P12").Value = 4.3
$ws.Range("R12").Value = 2.34
$ws.Range("T12").Value = 1.66
$ws.Range("Y12").Value = 20
$ws.Range("AC12").Value = 20
$ws.Range("AF12").Value = 140
$ws.Range("AK12").Value = 140
$ws.Range("AL12").Value = 90
$ws.Range("AM12").Value = 80
$ws.Range("AO12").Value = 2.8

# Row 13
$ws.Range("AB13").Value = 9.4
$ws.Range("AH13").Value = 70

# Row 14
$ws.Range("F14").Value = 1.45
$ws.Range("G14").Value = 1.47
$ws.Range("H14").Value = 9
$ws.Range("I14").Value = 9.4
$ws.Range("J14").Value = 4.7
$ws.Range("K14").Value = 4.9
$ws.Range("N14").Value = 3.95
$ws.Range("O14").Value = 1.31
$ws.Range("Q14").Value = 1.88
$ws.Range("T14").Value = 2.12
$ws.Range("U14").Value = 1.8
$ws.Range("V14").Value = 1.11
$ws.Range("W14").Value = 3.15
$ws.Range("Z14").Value = 80
$ws.Range("AB14").Value = 7.6
$ws.Range("AD14").Value = 34
$ws.Range("AE14").Value = 160
$ws.Range("AH14").Value = 30
$ws.Range("AL14").Value = 970
$ws.Range("AN14").Value = 7.8
$ws.Range("AO14").Value = 230

# Row 15
$ws.Range("F15").Value = 7.2
$ws.Range("G15").Value = 7.6
$ws.Range("H15").Value = 1.49
$ws.Range("J15").Value = 4.9
$ws.Range("N15").Value = 5.1
$ws.Range("O15").Value = 1.21
$ws.Range("R15").Value = 1.56
$ws.Range("S15").Value = 2.58
$ws.Range("T15").Value = 1.79
$ws.Range("U15").Value = 2.12
$ws.Range("AD15").Value = 9.800000000000001
$ws.Range("AI15").Value = 29
$ws.Range("AK15").Value = 100
$ws.Range("AN15").Value = 90
$ws.Range("AO15").Value = 6.4

# Row 16
$ws.Range("F16").Value = 1.39
$ws.Range("G16").Value = 1.45
$ws.Range("H16").Value = 8.4
$ws.Range("I16").Value = 11.5
$ws.Range("J16").Value = 5.2
$ws.Range("K16").Value = 5.9
$ws.Range("N16").Value = 2.38
$ws.Range("P16").Value = 2.38
$ws.Range("Q16").Value = 1.62
$ws.Range("R16").Value = 1.46
$ws.Range("S16").Value = 2.32
$ws.Range("T16").Value = 1.9
$ws.Range("U16").Value = 1.94
$ws.Range("V16").Value = 1.09
$ws.Range("W16").Value = 3.2

# Row 17
$ws.Range("N17").Value = 3.75
$ws.Range("R17").Value = 1.37

# Row 18
$ws.Range("F18").Value = 1.41
$ws.Range("G18").Value = 1.47
$ws.Range("H18").Value = 8
$ws.Range("K18").Value = 5.8
$ws.Range("O18").Value = 1.21
$ws.Range("P18").Value = 2.36
$ws.Range("U18").Value = 1.96
$ws.Range("W18").Value = 3.1
$ws.Range("AA18").Value = 300
$ws.Range("AE18").Value = 140
$ws.Range("AM18").Value = 140
$ws.Range("AO18").Value = 160
